$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.734.39"
$ws.Range("E2").Value = "  -3.07%  "
$ws.Range("D3").Value = "2.577.05"
$ws.Range("E3").Value = "  -5.25%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "547.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.13"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.22%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.597"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.18%  "
$ws.Range("E9").Value = "  -2.98%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.162"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.74%  "
$ws.Range("E11").Value = "  -2.04%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.364"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.23%  "
$ws.Range("D13").Value = "3.027.77"
$ws.Range("E13").Value = "  -5.40%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.41"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.71%  "
$ws.Range("D15").Value = "61.576.75"
$ws.Range("E15").Value = "  -3.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000144"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.12%  "
$ws.Range("D17").Value = "2.574.47"
$ws.Range("E17").Value = "  -5.34%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.53"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.41%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.54"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "337.01"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.83%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.06"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.71%  "
$ws.Range("E22").Value = "  +0.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.494"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.72%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.21"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.56%  "
$ws.Range("E25").Value = "  -1.37%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.08"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.52"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.33%  "
$ws.Range("D29").Value = "0.0₃0839"
$ws.Range("E29").Value = "  -4.46%  "
$ws.Range("E30").Value = "  -2.47%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.29"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "160.76"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.42%  "
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.74"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "19.16"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.37%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.40"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.66%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.79"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.30%  "
$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "333.63"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.74%  "
$ws.Range("B39").Value = "SuiNetwork"
$ws.Range("C39").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.933"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.46%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.92"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.47%  "
$ws.Range("E41").Value = "  -0.88%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "37.48"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.26%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.57"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.90%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.998"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "2.127.70"
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.604"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.58%  "
$ws.Range("E47").Value = "  -1.31%  "
$ws.Range("E48").Value = "  -4.15%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.55"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.18%  "
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0965"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.75%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0239"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.67%  "
